$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.054.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4828"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3818"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07356"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9326"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07842"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.903.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.502"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.604"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.088.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.153"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.127.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.099"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.958"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08903"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.359"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7655"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.673"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.605"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.098"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05293"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5479"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.979"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.006"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.446"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4827"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.653"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06100"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
